$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5, shifting existing rows 5-79 down to 6-80
$ws.Rows("5").Insert()

# Populate the newly inserted row 5 with the new weekly data entry
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44630
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 100112030
$ws.Range("G5").Value = "Poroto granado"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 21000
$ws.Range("L5").Value = 23000
$ws.Range("M5").Value = 22000
$ws.Range("N5").Value = "$/malla 25 kilos"
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 880
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = "Hortaliza"
